$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy formatting from the (now-shifted) original columns F:G into new D:E so styles match per row
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D, E) with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 177500
$ws.Range("E8").Value = 155500
$ws.Range("D9").Value = 400
$ws.Range("E9").Value = 1200
$ws.Range("D10").Value = 177100
$ws.Range("E10").Value = 154300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 97400
$ws.Range("E17").Value = 126000
$ws.Range("D18").Value = 80100
$ws.Range("E18").Value = 29500
$ws.Range("D20").Value = 14300
$ws.Range("E20").Value = 13600
$ws.Range("D21").Value = 98300
$ws.Range("E21").Value = 46900
$ws.Range("D22").Value = 600
$ws.Range("E22").Value = 600
$ws.Range("D23").Value = 93800
$ws.Range("E23").Value = 42500
$ws.Range("D24").Value = 18300
$ws.Range("E24").Value = 8700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 75500
$ws.Range("E26").Value = 33800
$ws.Range("D27").Value = 75500
$ws.Range("E27").Value = 33800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -14300
$ws.Range("E32").Value = -13600
$ws.Range("D33").Value = 75500
$ws.Range("E33").Value = 33800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 75500
$ws.Range("E35").Value = 33800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 120300
$ws.Range("E41").Value = 44900
$ws.Range("D42").Value = 69000
$ws.Range("E42").Value = 68800
$ws.Range("D43").Value = 46800
$ws.Range("E43").Value = 66900
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 73300
$ws.Range("E45").Value = 74700
$ws.Range("D46").Value = 309500
$ws.Range("E46").Value = 255200
$ws.Range("D47").Value = 900100
$ws.Range("E47").Value = 882900
$ws.Range("D48").Value = 111000
$ws.Range("E48").Value = 110900
$ws.Range("D49").Value = 2900
$ws.Range("E49").Value = 2900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 500
$ws.Range("E52").Value = 700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1324000
$ws.Range("E54").Value = 1252700
$ws.Range("D57").Value = 14300
$ws.Range("E57").Value = 17000
$ws.Range("D58").Value = 36500
$ws.Range("E58").Value = 6600
$ws.Range("D59").Value = 29800
$ws.Range("E59").Value = 31400
$ws.Range("D60").Value = 80500
$ws.Range("E60").Value = 54900
$ws.Range("D61").Value = 23400
$ws.Range("E61").Value = 55000
$ws.Range("D62").Value = 6500
$ws.Range("E62").Value = 4400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 110400
$ws.Range("E66").Value = 114300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1082200
$ws.Range("E72").Value = 1006600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1213600
$ws.Range("E76").Value = 1138400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 75500
$ws.Range("E81").Value = 33800
$ws.Range("D83").Value = 3900
$ws.Range("E83").Value = 3800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 102500
$ws.Range("E89").Value = 40700
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -22500
$ws.Range("E94").Value = -148300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -4600
$ws.Range("E100").Value = -7500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 75500
$ws.Range("E102").Value = -115100

# Data corrections beyond the simple shift (rows 47 & 91)
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "NA"
$ws.Range("H47").Value = "NA"
$ws.Range("I47").Value = "NA"
$ws.Range("J47").Value = "NA"
$ws.Range("F91").Value = -200
$ws.Range("G91").Value = -200
$ws.Range("H91").Value = -200
$ws.Range("I91").Value = -800
$ws.Range("J91").Value = -1500
